$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.178.90"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "3.545.45"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.50%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("E14").Value = "  -2.36%  "
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("E19").Value = "  -4.12%  "
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "419.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("E34").Value = "  -3.71%  "
$ws.Range("E36").Value = "  -2.96%  "
$ws.Range("E37").Value = "  -10.89%  "
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("E39").Value = "  -5.35%  "
$ws.Range("E40").Value = "  -7.18%  "
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("E45").Value = "  -5.75%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("E47").Value = "  -3.91%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  -7.08%  "
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("E51").Value = "  -7.24%  "
